# The edit rotates the data of rows 17-21 on the "Artfynd" sheet:
#   new row17 <- old row19
#   new row18 <- old row17
#   new row19 <- old row18
#   new row20 <- old row21
#   new row21 <- old row20
# (two independent cycles: 17 -> 18 -> 19 -> 17, and 20 -> 21 -> 20)
#
# We use whole-row range copies (columns A:AY, which covers the sheet's
# used range) via .Value2, staging through scratch rows far below the
# used data (rows 200/201) so that a source row is never clobbered
# before it has been read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns Y and AA hold ISO-formatted date text ("2026-02-03"). Those two
# columns are identical across rows 17-21 both before and after the edit,
# so we simply skip them during the row copy (copying A:X, Z:Z and AB:AY)
# to avoid Excel's automatic text->date coercion when such a literal is
# pushed back through .Value2.
function Copy-Row($srcRow, $dstRow) {
    $src1 = $ws.Range("A" + $srcRow + ":X" + $srcRow)
    $dst1 = $ws.Range("A" + $dstRow + ":X" + $dstRow)
    $dst1.Value2 = $src1.Value2

    $src2 = $ws.Range("Z" + $srcRow)
    $dst2 = $ws.Range("Z" + $dstRow)
    $dst2.Value2 = $src2.Value2

    $src3 = $ws.Range("AB" + $srcRow + ":AY" + $srcRow)
    $dst3 = $ws.Range("AB" + $dstRow + ":AY" + $dstRow)
    $dst3.Value2 = $src3.Value2
}

# Cycle A: 17 -> 18 -> 19 -> 17
Copy-Row 17 200          # stash old row17
Copy-Row 19 17           # row17 <- old row19
Copy-Row 18 19           # row19 <- old row18
Copy-Row 200 18          # row18 <- old row17 (stashed)

# Cycle B: 20 -> 21 -> 20
Copy-Row 21 201           # stash old row21
Copy-Row 20 21            # row21 <- old row20
Copy-Row 201 20            # row20 <- old row21 (stashed)

# Clean up the scratch rows used for staging.
$ws.Range("A200:AY201").Clear()
